# Weekly update: insert a new price observation as row 208 on the
# "Terminal Hortofrutícola Agro Chillán - Zapallo italiano" sheet,
# pushing the previously-existing rows 208:231 down to 209:232.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 208 (shifts 208:231 -> 209:232).
$ws.Rows.Item(208).Insert()

# Populate the new row 208 with the latest weekly observation.
$ws.Cells.Item(208, 1).Value = 7
$ws.Cells.Item(208, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(208, 3).Value = "Ñuble"
$ws.Cells.Item(208, 4).Value = 44776
$ws.Cells.Item(208, 5).Value = 16
$ws.Cells.Item(208, 6).Value = 100112032
$ws.Cells.Item(208, 7).Value = "Zapallo italiano"
$ws.Cells.Item(208, 8).Value = "Sin especificar"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 60
$ws.Cells.Item(208, 11).Value = 19000
$ws.Cells.Item(208, 12).Value = 20000
$ws.Cells.Item(208, 13).Value = 19500
$ws.Cells.Item(208, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(208, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(208, 16).Value = 390
$ws.Cells.Item(208, 17).Value = 50
$ws.Cells.Item(208, 18).Value = "Hortaliza"
